# Generate Report for Handoff
# Adds a new tracked file (6f1aaa4b-...) to the localization-status report:
#   - a new row on "Overview"
#   - a new row on "zh-cn"
#   - a new row on "de-de"
# and widens a couple of "handoff datetime" columns to fit the new values.

$wb = $excel.ActiveWorkbook

$commitSha = "851b6b58f8dac8b54d3ed41df6f91e8d320d3d5c"
$newFile   = "6f1aaa4b-80b6-4d6c-98bd-d626e64a58aaooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newFileUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/$commitSha/e2e/$newFile"

# ---------------------------------------------------------------------------
# Overview sheet: new row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", "e2e\$newFile") | Out-Null
$wsOverview.Range("B3").Style = "HyperLink"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-13 04:32:37"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet: new row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $newFile
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newFileUrl, "", "", "e2e\$newFile") | Out-Null
$wsZh.Range("A3").Style = "HyperLink"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = "6f1aaa4b-80b6-4d6c-98bd-d626e64a58aaooooooooooooooooooooooooooooooooooooooooooo.e2a0969aff0947a7549445ca487721abacc98604.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-13 04:32:29"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M3").Value = "True"
$wsZh.Range("O3").Value = "False"

$wsZh.Columns.Item(3).ColumnWidth = 16.33

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet: new row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $newFile
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newFileUrl, "", "", "e2e\$newFile") | Out-Null
$wsDe.Range("A3").Style = "HyperLink"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = "6f1aaa4b-80b6-4d6c-98bd-d626e64a58aaooooooooooooooooooooooooooooooooooooooooooo.e2a0969aff0947a7549445ca487721abacc98604.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-13 04:32:37"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M3").Value = "True"
$wsDe.Range("O3").Value = "False"

$wsDe.Columns.Item(3).ColumnWidth = 16.33

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
